# Update "Training Dashboard" sheet: decrement PERIOD TO EXPIRE (col H)
# by 1 day and bump LAST UPDATE (col I) from 03-Nov-2025 to 04-Nov-2025
# for rows 3 through 28 (new progress as of 04-Nov-2025).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Force column I to be treated as plain text so the "dd-mmm-yyyy" style
# strings are not auto-converted into date serial numbers.
$ws.Range("I3:I28").NumberFormat = "@"

for ($row = 3; $row -le 28; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # column I - LAST UPDATE

    $hCell.Value2 = $hCell.Value2 - 1
    $iCell.Value2 = "04-Nov-2025"
}
